$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 7.134879666666667
$ws.Range("H2").Value = 21.404639
$ws.Range("I2").Value = 0.07716103050836744
$ws.Range("J2").Value = 0.07716103050836744
$ws.Range("M2").Value = 4.391406
$ws.Range("N2").Value = 8.782812
$ws.Range("O2").Value = 0.04903950466333907
$ws.Range("P2").Value = 0.03453691982955612
$ws.Range("Q2").Value = 31.332153377478
$ws.Range("R2").Value = 187.992920264868
$ws.Range("S2").Value = 0.003783938715443133
$ws.Range("T2").Value = 0.00266490432463342

# Row 3
$ws.Range("G3").Value = 7.134879666666667
$ws.Range("H3").Value = 21.404639
$ws.Range("I3").Value = 0.07716103050836744
$ws.Range("J3").Value = 0.07716103050836744
$ws.Range("M3").Value = 7.979832333333334
$ws.Range("O3").Value = 0.08911201217176548
$ws.Range("P3").Value = 0.09413801509686184
$ws.Range("Q3").Value = 56.93514345850923
$ws.Range("R3").Value = 512.416291126583
$ws.Range("S3").Value = 0.006875974689847607
$ws.Range("T3").Value = 0.007263786254886111

# Row 4
$ws.Range("G4").Value = 7.134879666666667
$ws.Range("H4").Value = 21.404639
$ws.Range("I4").Value = 0.07716103050836744
$ws.Range("J4").Value = 0.07716103050836744
$ws.Range("M4").Value = 10.281678
$ws.Range("N4").Value = 30.845034
$ws.Range("O4").Value = 0.114817075949696
$ws.Range("P4").Value = 0.1212928691173092
$ws.Range("Q4").Value = 73.35853530141399
$ws.Range("R4").Value = 660.226817712726
$ws.Range("S4").Value = 0.008859403900236035
$ws.Range("T4").Value = 0.009359082774408113

# Row 5
$ws.Range("G5").Value = 7.134879666666667
$ws.Range("H5").Value = 21.404639
$ws.Range("I5").Value = 0.07716103050836744
$ws.Range("J5").Value = 0.07716103050836744
$ws.Range("M5").Value = 9.951477499999999
$ws.Range("N5").Value = 19.902955
$ws.Range("O5").Value = 0.1111296762969226
$ws.Range("P5").Value = 0.07826499772581527
$ws.Range("Q5").Value = 71.00259446804083
$ws.Range("R5").Value = 426.015566808245
$ws.Range("S5").Value = 0.008574880343131846
$ws.Range("T5").Value = 0.00603900787725894

# Row 6
$ws.Range("G6").Value = 7.134879666666667
$ws.Range("H6").Value = 21.404639
$ws.Range("I6").Value = 0.07716103050836744
$ws.Range("J6").Value = 0.07716103050836744
$ws.Range("M6").Value = 52.560594
$ws.Range("N6").Value = 157.681782
$ws.Range("O6").Value = 0.586952218622207
$ws.Range("P6").Value = 0.6200568865091892
$ws.Range("Q6").Value = 375.013513398522
$ws.Range("R6").Value = 3375.121620586698
$ws.Range("S6").Value = 0.04528983804806207
$ws.Range("T6").Value = 0.04784422833685888

# Row 7
$ws.Range("G7").Value = 7.134879666666667
$ws.Range("H7").Value = 21.404639
$ws.Range("I7").Value = 0.07716103050836744
$ws.Range("J7").Value = 0.07716103050836744
$ws.Range("M7").Value = 4.383347333333333
$ws.Range("N7").Value = 13.150042
$ws.Range("O7").Value = 0.04894951229606984
$ws.Range("P7").Value = 0.05171031172126829
$ws.Range("Q7").Value = 31.27465576053755
$ws.Range("R7").Value = 281.471901844838
$ws.Range("S7").Value = 0.003776994811646752
$ws.Range("T7").Value = 0.003990020940321973

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 81.06813666666666
$ws.Range("H8").Value = 243.20441
$ws.Range("I8").Value = 0.8767212985829614
$ws.Range("J8").Value = 0.8767212985829616
$ws.Range("M8").Value = 4.391406
$ws.Range("N8").Value = 8.782812
$ws.Range("O8").Value = 0.04903950466333907
$ws.Range("P8").Value = 0.03453691982955612
$ws.Range("Q8").Value = 356.00310176682
$ws.Range("R8").Value = 2136.01861060092
$ws.Range("S8").Value = 0.04299397821030782
$ws.Range("T8").Value = 0.03027925320202408

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 81.06813666666666
$ws.Range("H9").Value = 243.20441
$ws.Range("I9").Value = 0.8767212985829614
$ws.Range("J9").Value = 0.8767212985829616
$ws.Range("M9").Value = 7.979832333333334
$ws.Range("O9").Value = 0.08911201217176548
$ws.Range("P9").Value = 0.09413801509686184
$ws.Range("Q9").Value = 646.9101381757523
$ws.Range("R9").Value = 5822.19124358177
$ws.Range("S9").Value = 0.07812639903057089
$ws.Range("T9").Value = 0.08253280284174315

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 81.06813666666666
$ws.Range("H10").Value = 243.20441
$ws.Range("I10").Value = 0.8767212985829614
$ws.Range("J10").Value = 0.8767212985829616
$ws.Range("M10").Value = 10.281678
$ws.Range("N10").Value = 30.845034
$ws.Range("O10").Value = 0.114817075949696
$ws.Range("P10").Value = 0.1212928691173092
$ws.Range("Q10").Value = 833.5164772666599
$ws.Range("R10").Value = 7501.648295399939
$ws.Range("S10").Value = 0.100662575926116
$ws.Range("T10").Value = 0.1063400417213805

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 81.06813666666666
$ws.Range("H11").Value = 243.20441
$ws.Range("I11").Value = 0.8767212985829614
$ws.Range("J11").Value = 0.8767212985829616
$ws.Range("M11").Value = 9.951477499999999
$ws.Range("N11").Value = 19.902955
$ws.Range("O11").Value = 0.1111296762969226
$ws.Range("P11").Value = 0.07826499772581527
$ws.Range("Q11").Value = 806.7477380052583
$ws.Range("R11").Value = 4840.48642803155
$ws.Range("S11").Value = 0.09742975411414216
$ws.Range("T11").Value = 0.0686165904397693

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 81.06813666666666
$ws.Range("H12").Value = 243.20441
$ws.Range("I12").Value = 0.8767212985829614
$ws.Range("J12").Value = 0.8767212985829616
$ws.Range("M12").Value = 52.560594
$ws.Range("N12").Value = 157.681782
$ws.Range("O12").Value = 0.586952218622207
$ws.Range("P12").Value = 0.6200568865091892
$ws.Range("Q12").Value = 4260.98941767318
$ws.Range("R12").Value = 38348.90475905862
$ws.Range("S12").Value = 0.5145935113166117
$ws.Range("T12").Value = 0.5436170787356444

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 81.06813666666666
$ws.Range("H13").Value = 243.20441
$ws.Range("I13").Value = 0.8767212985829614
$ws.Range("J13").Value = 0.8767212985829616
$ws.Range("M13").Value = 4.383347333333333
$ws.Range("N13").Value = 13.150042
$ws.Range("O13").Value = 0.04894951229606984
$ws.Range("P13").Value = 0.05171031172126829
$ws.Range("Q13").Value = 355.3498006761355
$ws.Range("R13").Value = 3198.14820608522
$ws.Range("S13").Value = 0.04291507998521299
$ws.Range("T13").Value = 0.04533553164240008

# Row 14
$ws.Range("G14").Value = 4.264381
$ws.Range("H14").Value = 12.793143
$ws.Range("I14").Value = 0.04611767090867112
$ws.Range("J14").Value = 0.04611767090867113
$ws.Range("M14").Value = 4.391406
$ws.Range("N14").Value = 8.782812
$ws.Range("O14").Value = 0.04903950466333907
$ws.Range("P14").Value = 0.03453691982955612
$ws.Range("Q14").Value = 18.726628309686
$ws.Range("R14").Value = 112.359769858116
$ws.Range("S14").Value = 0.002261587737588114
$ws.Range("T14").Value = 0.001592762302898628

# Row 15
$ws.Range("G15").Value = 4.264381
$ws.Range("H15").Value = 12.793143
$ws.Range("I15").Value = 0.04611767090867112
$ws.Range("J15").Value = 0.04611767090867113
$ws.Range("M15").Value = 7.979832333333334
$ws.Range("O15").Value = 0.08911201217176548
$ws.Range("P15").Value = 0.09413801509686184
$ws.Range("Q15").Value = 34.02904538545234
$ws.Range("R15").Value = 306.2614084690711
$ws.Range("S15").Value = 0.004109638451346976
$ws.Range("T15").Value = 0.004341426000232589

# Row 16
$ws.Range("G16").Value = 4.264381
$ws.Range("H16").Value = 12.793143
$ws.Range("I16").Value = 0.04611767090867112
$ws.Range("J16").Value = 0.04611767090867113
$ws.Range("M16").Value = 10.281678
$ws.Range("N16").Value = 30.845034
$ws.Range("O16").Value = 0.114817075949696
$ws.Range("P16").Value = 0.1212928691173092
$ws.Range("Q16").Value = 43.844992311318
$ws.Range("R16").Value = 394.604930801862
$ws.Range("S16").Value = 0.005295096123343978
$ws.Range("T16").Value = 0.005593744621520585

# Row 17
$ws.Range("G17").Value = 4.264381
$ws.Range("H17").Value = 12.793143
$ws.Range("I17").Value = 0.04611767090867112
$ws.Range("J17").Value = 0.04611767090867113
$ws.Range("M17").Value = 9.951477499999999
$ws.Range("N17").Value = 19.902955
$ws.Range("O17").Value = 0.1111296762969226
$ws.Range("P17").Value = 0.07826499772581527
$ws.Range("Q17").Value = 42.4368915729275
$ws.Range("R17").Value = 254.621349437565
$ws.Range("S17").Value = 0.005125041839648628
$ws.Range("T17").Value = 0.003609399408787043

# Row 18
$ws.Range("G18").Value = 4.264381
$ws.Range("H18").Value = 12.793143
$ws.Range("I18").Value = 0.04611767090867112
$ws.Range("J18").Value = 0.04611767090867113
$ws.Range("M18").Value = 52.560594
$ws.Range("N18").Value = 157.681782
$ws.Range("O18").Value = 0.586952218622207
$ws.Range("P18").Value = 0.6200568865091892
$ws.Range("Q18").Value = 224.138398402314
$ws.Range("R18").Value = 2017.245585620826
$ws.Range("S18").Value = 0.02706886925753333
$ws.Range("T18").Value = 0.02859557943668603

# Row 19
$ws.Range("G19").Value = 4.264381
$ws.Range("H19").Value = 12.793143
$ws.Range("I19").Value = 0.04611767090867112
$ws.Range("J19").Value = 0.04611767090867113
$ws.Range("M19").Value = 4.383347333333333
$ws.Range("N19").Value = 13.150042
$ws.Range("O19").Value = 0.04894951229606984
$ws.Range("P19").Value = 0.05171031172126829
$ws.Range("Q19").Value = 18.69226308466733
$ws.Range("R19").Value = 168.230367762006
$ws.Range("S19").Value = 0.0022574374992101
$ws.Range("T19").Value = 0.002384759138546251
